$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "FirstName"
$ws.Range("B1").Value = "LastName"
$ws.Range("D1").Value = "DoB"
$ws.Range("F1").Value = "PhoneNumber"
$ws.Range("H1").Value = "IsGraduated"
